# Updated cryptos list on Sun Oct 29 05:13:46 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto
# tracker sheet, plus a ranking swap between Chainlink and WrappedEther
# (rows 13/14) where WrappedEther moved ahead of Chainlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal text (the Price
# column holds strings like "34.095.85" / "0.0686" which Excel would
# otherwise happily reinterpret as a number/date on assignment). Forcing
# the cell to Text format before the write keeps the exact string, then
# resetting the style back to Normal avoids leaving a stray text format
# behind on the cell.
function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$Value)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Value
    $range.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue $ws "D2" "34.095.85"
$ws.Range("E2").Value = "  -0.07%  "

# --- Row 3: Ethereum ---
Set-TextValue $ws "D3" "1.781.36"
$ws.Range("E3").Value = "  -0.53%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.27%  "

# --- Row 5: BNB ---
Set-TextValue $ws "D5" "225.34"
$ws.Range("E5").Value = "  -0.97%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  -0.08%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.28%  "

# --- Row 8: Solana ---
Set-TextValue $ws "D8" "31.94"
$ws.Range("E8").Value = "  -1.27%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  -1.72%  "

# --- Row 10: Dogecoin ---
Set-TextValue $ws "D10" "0.0686"
$ws.Range("E10").Value = "  -0.19%  "

# --- Row 11: TRON ---
Set-TextValue $ws "D11" "0.0949"
$ws.Range("E11").Value = "  +0.55%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
Set-TextValue $ws "D12" "2.038.10"
$ws.Range("E12").Value = "  -0.49%  "

# --- Row 13: was Chainlink, now WrappedEther ---
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D13" "1.792.54"
$ws.Range("E13").Value = "  +0.06%  "

# --- Row 14: was WrappedEther, now Chainlink ---
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D14" "10.90"
$ws.Range("E14").Value = "  -5.87%  "

# --- Row 15: Polygon ---
$ws.Range("E15").Value = "  -0.28%  "

# --- Row 16: WrappedBTC ---
Set-TextValue $ws "D16" "34.079.16"
$ws.Range("E16").Value = "  -0.06%  "

# --- Row 17: Polkadot ---
$ws.Range("E17").Value = "  -0.32%  "

# --- Row 18: Litecoin ---
Set-TextValue $ws "D18" "67.53"
$ws.Range("E18").Value = "  -0.96%  "

# --- Row 19: BitcoinCash ---
Set-TextValue $ws "D19" "245.58"
$ws.Range("E19").Value = "  +0.68%  "

# --- Row 20: ShibaInu ---
Set-TextValue $ws "D20" "0.0₃0787"
$ws.Range("E20").Value = "  +1.30%  "

# --- Row 21: Dai ---
$ws.Range("E21").Value = "  +0.29%  "

# --- Row 22: Avalanche ---
Set-TextValue $ws "D22" "10.85"
$ws.Range("E22").Value = "  +0.47%  "

# --- Row 23: Uniswap ---
$ws.Range("E23").Value = "  -0.32%  "

# --- Row 24: Toncoin ---
$ws.Range("E24").Value = "  -1.15%  "

# --- Row 25: Monero ---
$ws.Range("E25").Value = "  +0.05%  "

# --- Row 26: Cosmos ---
Set-TextValue $ws "D26" "7.12"
$ws.Range("E26").Value = "  -0.75%  "

# --- Row 27: EthereumClassic ---
$ws.Range("E27").Value = "  -0.14%  "

# --- Row 28: Stellar ---
$ws.Range("E28").Value = "  +0.12%  "

# --- Row 29: BinanceUSD ---
$ws.Range("E29").Value = "  +0.38%  "

# --- Row 30: PancakeSwap ---
$ws.Range("E30").Value = "  -1.39%  "

# --- Row 31: Hedera ---
Set-TextValue $ws "D31" "0.0516"
$ws.Range("E31").Value = "  -0.29%  "

# --- Row 32: Filecoin ---
$ws.Range("E32").Value = "  +1.18%  "

# --- Row 33: InternetComputer(DFINITY) ---
$ws.Range("E33").Value = "  +2.00%  "

# --- Row 34: LidoDAOToken ---
$ws.Range("E34").Value = "  -2.57%  "

# --- Row 35: Maker ---
Set-TextValue $ws "D35" "1.448.26"
$ws.Range("E35").Value = "  +2.49%  "

# --- Row 36: RenderToken ---
$ws.Range("E36").Value = "  +4.56%  "

# --- Row 37: ImmutableX ---
Set-TextValue $ws "D37" "0.650"
$ws.Range("E37").Value = "  -0.85%  "

# --- Row 38: VeChain ---
$ws.Range("E38").Value = "  +0.50%  "

# --- Row 39: TrustWalletToken ---
$ws.Range("E39").Value = "  -0.95%  "

# --- Row 40: HuobiToken ---
$ws.Range("E40").Value = "  +1.34%  "

# --- Row 41: Aave ---
Set-TextValue $ws "D41" "80.58"
$ws.Range("E41").Value = "  +0.41%  "

# --- Row 42: MXToken ---
$ws.Range("E42").Value = "  +1.29%  "

# --- Row 43: ARBITRUM ---
Set-TextValue $ws "D43" "0.913"
$ws.Range("E43").Value = "  -1.35%  "

# --- Row 44: InjectiveProtocol ---
Set-TextValue $ws "D44" "13.64"
$ws.Range("E44").Value = "  +2.38%  "

# --- Row 45: Kaspa ---
Set-TextValue $ws "D45" "0.0519"
$ws.Range("E45").Value = "  +2.10%  "

# --- Row 46: FraxShare ---
Set-TextValue $ws "D46" "6.06"
$ws.Range("E46").Value = "  +0.08%  "

# --- Row 47: WEMIXToken ---
$ws.Range("E47").Value = "  +0.29%  "

# --- Row 48: RocketPoolETH ---
Set-TextValue $ws "D48" "1.938.63"
$ws.Range("E48").Value = "  -0.57%  "

# --- Row 49: BabyDogeCoin ---
$ws.Range("E49").Value = "  -6.24%  "

# --- Row 50: Quant ---
$ws.Range("E50").Value = "  -2.90%  "

# --- Row 51: PaxDollar ---
$ws.Range("E51").Value = "  +0.28%  "
